$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Log")
$ws2 = $wb.Worksheets.Item("To Do")

# Sheet1 (Log): log a new completed entry describing the streamlit slider -> input field change
$ws1.Range("A5").Value = 44647
$ws1.Range("A5").NumberFormat = "d-mmm"
$ws1.Range("B5").Value = "Changed some streamlit sliders to input fields"

# Sheet2 (To Do): the "Check what is wrong..." item is done (moved to Log above, reworded),
# so shift the remaining todo up and add a new todo item at the top
$ws2.Range("A1").Value = "Investigate how to depliy streamlit prototype to AWS"
$ws2.Range("A2").Value = "Think what to do with owner edrpous"
$ws2.Range("A3").ClearContents()

# Restore selections to match the saved view state
$ws1.Activate()
$ws1.Range("B8").Select()
$ws2.Activate()
$ws2.Range("A5").Select()
